# Re-point the three "Table_0" tables (custom tableStyleId
# {3DA7F5AD-EB32-4D34-83E3-F521C89693B4}) at the built-in PowerPoint
# table style {5A4CF338-0256-4E8C-9F09-137C2797A291} (Medium Style 2 -
# Accent 1), on slides 14, 15 and 16.

$p = $ppt.ActivePresentation
$newStyleId = "{5A4CF338-0256-4E8C-9F09-137C2797A291}"
$targetSlides = 14, 15, 16

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
